# Mechanical Report edits
# Applies the set of wording/content changes described in the commit
# "Properly Finished Mech Report".

$d = $word.ActiveDocument

function FindReplace($searchText, $replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "WARNING: could not find [$searchText]"
    }
    return $ok
}

# 1) "...for an operating point that is used to run a simulation..."
#    -> "...for an operating point, which is in turn used to run a simulation..."
FindReplace `
    "for an operating point that is used to run a simulation of the quality control procedure" `
    "for an operating point, which is in turn used to run a simulation of the quality control procedure"

# 2) "...a linear gripping mechanism, providing 3.5 degrees of freedom."
#    -> "...a linear gripping mechanism, totaling to 3.5 degrees of freedom."
FindReplace "providing 3.5 degrees of freedom." "totaling to 3.5 degrees of freedom."

# 3) Add a new sentence referencing Figure 2 at the end of the 3D-Design
#    paragraph that talks about arm segment / claw dimensions.
$r = $d.Content
$found = $r.Find.Execute("approximate size of each marshmallow.")
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter(" This can be seen in ")
    $insPos = $r.End
    $figRange = $d.Range($insPos, $insPos)
    $figRange.InsertAfter("Figure 2")
    $figTextRange = $d.Range($insPos, $insPos + 8)
    $figTextRange.Font.Italic = $true
    $figTextRange.Font.ItalicBi = $true
    $endRange = $d.Range($figTextRange.End, $figTextRange.End)
    $endRange.InsertAfter(".")
} else {
    Write-Output "WARNING: marshmallow sentence not found"
}

# 4) "The custom parts can be 3D printed from ABS plastic. These parts
#    include the first and second arm segment, the gripper base..."
#    -> "The remaining custom parts can be 3D printed from ABS plastic.
#    These parts include the first and second arm segments, the gripper
#    base..."
FindReplace "The custom parts can be 3D printed from ABS plastic." `
            "The remaining custom parts can be 3D printed from ABS plastic."
FindReplace "These parts include the first and second arm segment, the gripper base" `
            "These parts include the first and second arm segments, the gripper base"

# 5) Mark the run that hosts the "Custom ABS Parts" figure (Figure 2,
#    the first InlineShape) as NoProof.
$shp = $d.InlineShapes.Item(1)
$shp.Range.NoProofing = $true

# 6) "The moment of inertia for each of the three arm motors is found by
#    removing the preceding components, placing the axis of rotation..."
#    -> "The moment of inertia for each of the three motors that control
#    the arm movement, is found by removing the preceding components,
#    then placing the axis of rotation..."
FindReplace `
    "The moment of inertia for each of the three arm motors is found by removing the preceding components, placing the axis of rotation at the origin," `
    "The moment of inertia for each of the three motors that control the arm movement, is found by removing the preceding components, then placing the axis of rotation at the origin,"

# 7) "To prevent any failure of the robot's components, the likely points
#    of stress undergo a simulated test."
#    -> "To prevent any failure of the robot's components in practice, the
#    likely points of stress undergo a simulated test."
FindReplace `
    "To prevent any failure of the robot’s components, the likely points of stress undergo a simulated test." `
    "To prevent any failure of the robot’s components in practice, the likely points of stress undergo a simulated test."

# 8) "...whereas a SimulationX model is. The simulation time..."
#    -> "...whereas a SimulationX model will suffice. The simulation time..."
FindReplace " model is. The simulation time" " model will suffice. The simulation time"

# 9) "All of the motors are combinations of cylinders and joint blocks."
#    -> "All of the motors are represented by combinations of cylinders
#    and joint blocks."
FindReplace "the motors are combinations of cylinders and joint blocks" `
            "the motors are represented by combinations of cylinders and joint blocks"

# 10) Remove the stray grammar-check markers that wrapped "In order to" by
#     deleting the whole paragraph (proofErr tags included) and retyping
#     its text fresh.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "In order to greatly limit simulation time*") {
        $pr = $p.Range
        $insertPos = $pr.Start
        $full = $d.Range($pr.Start, $pr.End)
        $full.Delete()
        $newR = $d.Range($insertPos, $insertPos)
        $newR.InsertBefore("In order to greatly limit simulation time, the gripping mechanism is reduced to a single rotating component that is directly attached to the motor’s shaft. The claws are held in a fixed position. Since the gripping mechanism requires no feedback control, but a more primitive on/off control, this omission does not affect the accuracy of the control model.`r")
        break
    }
}

Write-Output "Done."
